# difot_score.xlsx template export — replace the old ad-hoc date/amount/name
# sample rows with the DIFOT score header row (SAP No / Supplier name /
# DIFOT score / Month/Year) and drop the now-unused date-formatted sample
# data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two sample data rows (old rows 2 and 3) — only the header row
# remains afterwards.
$ws.Rows("2:3").Delete()

# Overwrite the header row with the new DIFOT score column titles.
$ws.Range("A1").Value = "SAP No"
$ws.Range("B1").Value = "Supplier name"
$ws.Range("C1").Value = "DIFOT score"
$ws.Range("D1").Value = "Month/Year"

# Settle the selection back on the top-left cell now that the old B2
# selection target no longer exists.
$ws.Range("A1").Select()

# Widen the saved window size (author resized/maximized the window before
# exporting the template).
$win = $excel.ActiveWindow
$win.Width = 51200
$win.Height = 28260

$wb.Save()
